$wb = $excel.ActiveWorkbook

# --- Sheet "JatHarcos" (sheet1): add new row 36 ---
$wsHarcos = $wb.Worksheets.Item("JatHarcos")
$wsHarcos.Range("A36").Value = "'tzuiopő"
$wsHarcos.Range("B36").Value = "'Harcos"
$wsHarcos.Range("C36").Value = "'katana"
$wsHarcos.Range("D36").Value = "'/Images/Karakterek/harcos1.png"

# --- Sheet "JatMagus" (sheet2): add new row 33 ---
$wsMagus = $wb.Worksheets.Item("JatMagus")
$wsMagus.Range("A33").Value = "'Uiorepzr"
$wsMagus.Range("B33").Value = "'Mágus"
$wsMagus.Range("C33").Value = "'szél botja"
$wsMagus.Range("D33").Value = "'/Images/Karakterek/magus1.png"

Write-Host "Rows added"
